$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute could not find text: $old"
    }
}

# 1. Key Qualifications bullet
Replace-Text "Over 15 years building Web applications" `
             "Over 15 years experience building Web applications"

# 2. Led the development ... (split run before hyperlink)
Replace-Text "Led the development of an interactive D3.js powered visualization tool that powered our" `
             "Led development of interactive"

# 3. hyperlink text "healthcare decision tool" -- use Hyperlink.TextToDisplay
#    so the Link character style (w:rStyle) on the hyperlink run is preserved
#    instead of being dropped by a plain Find/Replace.
$foundHyperlink = $false
foreach ($h in $d.Hyperlinks) {
    if ($h.TextToDisplay -eq "healthcare decision tool") {
        $h.TextToDisplay = "visual healthcare decision support tool"
        $foundHyperlink = $true
    }
}
if (-not $foundHyperlink) {
    throw "Could not find hyperlink with text: healthcare decision tool"
}

# 4. Javascript framework bullet
Replace-Text "Developed a foundational Javascript framework for external contractors (built on on-top of Backbone.js) to ensure consistent approach to frontend javascript applications" `
             "Created Javascript framework (built on on top of Backbone.js) to ensure consistent approach to frontend Javascript applications"

# 5. automated frontend testing suite bullet
Replace-Text "Established automated frontend testing suite for our growing single page applications using Mocha/Chai/Karma" `
             "Established automated frontend testing suite using Mocha/Chai/Karma"

# 6. Spearheaded open source libraries bullet
Replace-Text "Spearheaded the development and Open Sourcing of several libraries used to grow our Backbone.js testing capabilities" `
             "Spearheaded the development of several open source libraries related to expanding Backbone.js testing capabilities"

# 7. Architected core time-series bullet
Replace-Text "Architected the core time-series data system to replace an aging MySQL solution with a highly scalable service backed by Clojure, Cassandra and RabbitMQ" `
             "Architected core time-series data system to replace an aging MySQL solution with a highly scalable service backed by Clojure, Cassandra and RabbitMQ"

# 8. Established Agile process bullet
Replace-Text "Established an Agile process to facilitate a rapid development and release cycle that significantly increased customer and company responsiveness and transition us from a release every few months, to multiple releases a week" `
             "Established Agile process to facilitate rapid development and release cycle from a release every few months to multiple releases a week"

# 9. Successfully designed, developed... merger bullet
Replace-Text "Successfully designed, developed, and led the technical merger of our product with our acquired company over the course of two months in order to differentiate us at our largest industry trade show in an increasingly dense competitive landscape" `
             "Successfully designed, developed, and led the technical merger of our product with an acquired company"

# 10. Oversaw and negotiated bullet
Replace-Text "Oversaw and negotiated a software staff increase as well as a company merger from five software developers to a fifteen member software development team consisting of QA, Design, DevOps, and Developers" `
             "Oversaw and negotiated software staff increase as well as company merger from five software developers to a fifteen member software development team consisting of QA, Design, DevOps, and Developers"

# 11. Handled budgeting bullet
Replace-Text "Handled budgeting, project scheduling and application & server architecture" `
             "Handled budgeting, project scheduling, application and server architecture"

# 12. Split "Demonstrated leadership..." bullet into two separate bullets
Replace-Text "Demonstrated leadership as technical lead; managing multiple project timelines and deliverables; working with colleagues to identify and develop solutions to operational requirements and shortcomings" `
             ("Managed multiple project timelines and deliverables`r" + "Worked with colleagues to identify and develop solutions to operational requirements and shortcomings")

# 13. Initiated transition bullet
Replace-Text "Initiated transition to open source platform" `
             "Initiated transition to open source course platform"

# 14. Designed and provided ongoing development (Internet site) bullet
Replace-Text "Designed and provided ongoing development and maintenance of national award-winning dynamic Internet site to allow customers to access company resources (JSP with DB2 backend)" `
             "Designed and provided ongoing development and maintenance of national award-winning dynamic Web site"

# 15. Created and introduced Intranet site bullet
Replace-Text "Created and introduced, with continuing development and maintenance, a unique Intranet site providing executives and employees access to critical business and process information. (ASP with Oracle backend)" `
             "Created a unique Intranet site providing executives and employees access to critical business and process information"

# 16. Developed numerous custom applications bullet
Replace-Text "Developed numerous custom applications and tools to solve unique client problems involving the conversion, automation, and filtering of business data (Visual Basic, PHP, ASP, Excel, Access)" `
             "Developed numerous custom applications and tools to solve unique client problems involving the conversion, automation, and filtering of business data"
